$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '44.277.88'
$ws.Cells.Item(2, 5).Value = '  -0.14%  '
$ws.Cells.Item(3, 4).Value = '2.263.89'
$ws.Cells.Item(3, 5).Value = '  -0.33%  '
$ws.Cells.Item(4, 5).Value = '  +0.11%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '317.29'
$cell.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.79%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '100.09'
$cell.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -4.84%  '
$ws.Cells.Item(7, 5).Value = '  -2.09%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  -5.44%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '36.36'
$cell.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -6.08%  '
$ws.Cells.Item(11, 5).Value = '  -1.74%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.44'
$cell.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -5.60%  '
$ws.Cells.Item(13, 5).Value = '  -2.38%  '
$ws.Cells.Item(14, 4).Value = '2.608.25'
$ws.Cells.Item(14, 5).Value = '  -0.40%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.851'
$cell.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -3.71%  '
$ws.Cells.Item(16, 4).Value = '2.257.18'
$ws.Cells.Item(16, 5).Value = '  -0.81%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '14.06'
$cell.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -3.56%  '
$ws.Cells.Item(18, 4).Value = '44.214.48'
$ws.Cells.Item(18, 5).Value = '  -0.08%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '13.22'
$cell.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -4.73%  '
$ws.Cells.Item(20, 5).Value = '  -2.04%  '
$ws.Cells.Item(21, 5).Value = '  -1.69%  '
$ws.Cells.Item(22, 5).Value = '  -0.78%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '239.82'
$cell.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.31%  '
$ws.Cells.Item(24, 5).Value = '  -6.06%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.05'
$cell.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -7.81%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.01'
$cell.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +0.54%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.29'
$cell.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '38.90'
$cell.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.18%  '
$ws.Cells.Item(29, 5).Value = '  -3.98%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.13'
$cell.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -5.54%  '
$ws.Cells.Item(31, 5).Value = '  -2.25%  '
$ws.Cells.Item(32, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.53'
$cell.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +12.93%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0848'
$cell.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -3.87%  '
$ws.Cells.Item(34, 2).Value = 'Monero'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '155.18'
$cell.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -4.86%  '
$ws.Cells.Item(35, 5).Value = '  -3.31%  '
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.93'
$cell.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -4.06%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.110'
$cell.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -5.46%  '
$ws.Cells.Item(38, 5).Value = '  -2.17%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '15.52'
$cell.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.44%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.54'
$cell.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -9.89%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.99'
$cell.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -9.36%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0311'
$cell.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -5.63%  '
$ws.Cells.Item(43, 5).Value = '  +0.09%  '
$ws.Cells.Item(44, 4).Value = '1.734.39'
$ws.Cells.Item(44, 5).Value = '  -2.68%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '84.03'
$cell.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -2.58%  '
$ws.Cells.Item(46, 5).Value = '  -5.53%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.22'
$cell.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -3.79%  '
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '15.39'
$cell.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +3.80%  '
$ws.Cells.Item(49, 2).Value = 'Aave'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '102.69'
$cell.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -1.44%  '
$ws.Cells.Item(50, 2).Value = 'MultiversX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '57.27'
$cell.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -5.61%  '
$ws.Cells.Item(51, 2).Value = 'ordi'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '71.71'
$cell.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -5.18%  '
